# Auto-generated edit script to update cryptos worksheet values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.576.02'
$ws.Range("E2").Value = '  +4.16%  '

# Row 3
$ws.Range("D3").Value = '1.846.06'
$ws.Range("E3").Value = '  +3.21%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.031'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +2.82%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '319.21'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.38%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.028'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.47%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4371'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +2.55%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3743'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +3.28%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07396'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.39%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8760'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.80%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '21.50'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.71%  '

# Row 12
$ws.Range("D12").Value = '1.849.58'
$ws.Range("E12").Value = '  +0.68%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.491'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.32%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.693'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.94%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.07153'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +3.93%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '82.75'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.81%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.033'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.99%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000009008'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.71%  '

# Row 19
$ws.Range("E19").Value = '  +2.47%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '15.41'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.74%  '

# Row 21
$ws.Range("D21").Value = '27.585.98'
$ws.Range("E21").Value = '  +4.09%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.268'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.52%  '

# Row 23
$ws.Range("E23").Value = '  +1.57%  '

# Row 24
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.064.52'
$ws.Range("E24").Value = '  +0.49%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '157.85'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.82%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.928'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +5.88%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.72'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.16%  '

# Row 28
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.257'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.15%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.939'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.22%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '116.22'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.34%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09087'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.21%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7683'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +3.16%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.204'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +5.39%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.500'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.52%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.879'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.29%  '

# Row 36
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.029'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.62%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.145'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.53%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01976'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.07%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05261'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.08%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5175'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.90%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.789'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +7.18%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.1673'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.27%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '6.667'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +4.11%  '

# Row 44
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '8.540'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.66%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '108.85'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.13%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '10.56'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.89%  '

# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.714'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.67%  '

# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.4650'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.90%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.06368'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.66%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.897'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +7.71%  '

# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '39.52'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.88%  '
